# Auto-generated edit script applying the Excalibur_Profits.xlsx diff
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) on specific rows
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1552.5454
$ws.Range("I98").Value = 1552.5454
$ws.Range("K98").Value = 1552.5454
$ws.Range("M98").Value = -54.54539999999997
# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1552.5454
$ws.Range("I122").Value = 1552.5454
$ws.Range("K122").Value = 4657.6362
$ws.Range("M122").Value = -2207.6362
# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("N125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 6107.44
$ws.Range("I132").Value = 6093.298
$ws.Range("J132").Value = 6329
$ws.Range("K132").Value = 18279.894
$ws.Range("L132").Value = 18987
$ws.Range("M132").Value = -15749.894
$ws.Range("N132").Value = -24047

$ws = $wb.Worksheets.Item("ARM")
# Row 26 (Leve Item ID 2241)
$ws.Range("H26").Value = 229.14285
$ws.Range("I26").Value = 229.14285
$ws.Range("K26").Value = 229.14285
$ws.Range("M26").Value = 100.85715
# Row 36 (Leve Item ID 3068)
$ws.Range("H36").Value = 2755.3333
$ws.Range("I36").Value = 1471.1428
$ws.Range("J36").Value = 7250
$ws.Range("K36").Value = 1471.1428
$ws.Range("L36").Value = 7250
$ws.Range("M36").Value = -1125.1428
$ws.Range("N36").Value = -7942
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 4260.143
$ws.Range("I45").Value = 4529.6313
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 4529.6313
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -4152.6313
$ws.Range("N45").Value = -2454
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4937
$ws.Range("I61").Value = 4375
$ws.Range("J61").Value = 5499
$ws.Range("K61").Value = 4375
$ws.Range("L61").Value = 5499
$ws.Range("M61").Value = -4163
$ws.Range("N61").Value = -5923
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4509.643
$ws.Range("I132").Value = 4043.6
$ws.Range("K132").Value = 12130.8
$ws.Range("M132").Value = -9600.799999999999
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4937
$ws.Range("I136").Value = 4375
$ws.Range("J136").Value = 5499
$ws.Range("K136").Value = 13125
$ws.Range("L136").Value = 16497
$ws.Range("M136").Value = -10575
$ws.Range("N136").Value = -21597

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 413.90475
$ws.Range("I22").Value = 436.3684
$ws.Range("J22").Value = 200.5
$ws.Range("K22").Value = 436.3684
$ws.Range("L22").Value = 200.5
$ws.Range("M22").Value = -263.3684
$ws.Range("N22").Value = -546.5
# Row 106 (Leve Item ID 18664)
$ws.Range("H106").Value = 23171
$ws.Range("J106").Value = 23171
$ws.Range("L106").Value = 23171
$ws.Range("N106").Value = -25695
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3589.45
$ws.Range("I107").Value = 3592.9285
$ws.Range("J107").Value = 3581.3333
$ws.Range("K107").Value = 3592.9285
$ws.Range("L107").Value = 3581.3333
$ws.Range("M107").Value = -1672.9285
$ws.Range("N107").Value = -7421.3333
# Row 112 (Leve Item ID 25788)
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("N112").Value = 0
$ws.Range("L112").ClearContents()
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 15553.177
$ws.Range("I134").Value = 17441.715
$ws.Range("J134").Value = 6740
$ws.Range("K134").Value = 52325.145
$ws.Range("L134").Value = 20220
$ws.Range("M134").Value = -49790.145
$ws.Range("N134").Value = -25290
# Row 135 (Leve Item ID 41992)
$ws.Range("H135").Value = 97755.8
$ws.Range("J135").Value = 97755.8
$ws.Range("L135").Value = 97755.8
$ws.Range("N135").Value = -107895.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 2583.625
$ws.Range("I16").Value = 2411.5
$ws.Range("K16").Value = 2411.5
$ws.Range("M16").Value = -2124.5
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 6905.4155
$ws.Range("I31").Value = 953.7818
$ws.Range("J31").Value = 21784.5
$ws.Range("K31").Value = 953.7818
$ws.Range("L31").Value = 21784.5
$ws.Range("M31").Value = -658.7818
$ws.Range("N31").Value = -22374.5
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 6905.4155
$ws.Range("I34").Value = 953.7818
$ws.Range("J34").Value = 21784.5
$ws.Range("K34").Value = 953.7818
$ws.Range("L34").Value = 21784.5
$ws.Range("M34").Value = -751.7818
$ws.Range("N34").Value = -22188.5
# Row 38 (Leve Item ID 1637)
$ws.Range("H38").Value = 5657.5
$ws.Range("I38").Value = 1900
$ws.Range("K38").Value = 1900
$ws.Range("M38").Value = -1523
# Row 46 (Leve Item ID 1637)
$ws.Range("H46").Value = 5657.5
$ws.Range("I46").Value = 1900
$ws.Range("K46").Value = 1900
$ws.Range("M46").Value = -1689
# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 5467.273
$ws.Range("I86").Value = 4238.8
$ws.Range("K86").Value = 4238.8
$ws.Range("M86").Value = -3115.8
# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 5467.273
$ws.Range("I89").Value = 4238.8
$ws.Range("K89").Value = 21194
$ws.Range("M89").Value = -15578
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 2583.625
$ws.Range("I113").Value = 2411.5
$ws.Range("K113").Value = 2411.5
$ws.Range("M113").Value = -241.5

$ws = $wb.Worksheets.Item("CUL")
# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 1539.9286
$ws.Range("J129").Value = 2321.8
$ws.Range("L129").Value = 6965.400000000001
$ws.Range("N129").Value = -16965.4

$ws = $wb.Worksheets.Item("GSM")
# Row 63 (Leve Item ID 11048)
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("N63").Value = 0
$ws.Range("L63").ClearContents()
# Row 66 (Leve Item ID 11048)
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("L66").ClearContents()
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 74685.22
$ws.Range("I80").Value = 87211.164
$ws.Range("K80").Value = 87211.164
$ws.Range("M80").Value = -86213.164
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 74685.22
$ws.Range("I83").Value = 87211.164
$ws.Range("K83").Value = 436055.82
$ws.Range("M83").Value = -431063.82
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 7193.9414
$ws.Range("I132").Value = 7393.5625
$ws.Range("K132").Value = 22180.6875
$ws.Range("M132").Value = -19650.6875

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 24523.8
$ws.Range("I22").Value = 59082
$ws.Range("J22").Value = 1485
$ws.Range("K22").Value = 59082
$ws.Range("L22").Value = 1485
$ws.Range("M22").Value = -58787
$ws.Range("N22").Value = -2075
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 24523.8
$ws.Range("I27").Value = 59082
$ws.Range("J27").Value = 1485
$ws.Range("K27").Value = 59082
$ws.Range("L27").Value = 1485
$ws.Range("M27").Value = -58975
$ws.Range("N27").Value = -1699
# Row 101 (Leve Item ID 18549)
$ws.Range("H101").Value = 25314.285
$ws.Range("J101").Value = 25314.285
$ws.Range("L101").Value = 25314.285
$ws.Range("N101").Value = -31804.285
# Row 110 (Leve Item ID 25809)
$ws.Range("H110").Value = 66506.73
$ws.Range("J110").Value = 66506.73
$ws.Range("L110").Value = 66506.73
$ws.Range("N110").Value = -74686.73
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 83769.53999999999
$ws.Range("J122").Value = 257501
$ws.Range("L122").Value = 772503
$ws.Range("N122").Value = -777403
# Row 128 (Leve Item ID 34582)
$ws.Range("H128").Value = 96952.664
$ws.Range("I128").Value = 90000
$ws.Range("J128").Value = 98343.2
$ws.Range("K128").Value = 90000
$ws.Range("L128").Value = 98343.2
$ws.Range("M128").Value = -85020
$ws.Range("N128").Value = -108303.2

$ws = $wb.Worksheets.Item("WVR")
# Row 64 (Leve Item ID 11036)
$ws.Range("H64").Value = 44999
$ws.Range("J64").Value = 44999
$ws.Range("L64").Value = 44999
$ws.Range("N64").Value = -45495
# Row 67 (Leve Item ID 11036)
$ws.Range("H67").Value = 44999
$ws.Range("J67").Value = 44999
$ws.Range("L67").Value = 44999
$ws.Range("N67").Value = -46715
# Row 137 (Leve Item ID 42184)
$ws.Range("H137").Value = 99000
$ws.Range("J137").Value = 99000
$ws.Range("L137").Value = 99000
$ws.Range("N137").Value = -109200
